$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value, whether the text looks like a
# plain number (needs to be forced to Text format so Excel does not silently
# convert it to a numeric value / lose trailing zeros).
$updates = @(
    @{ Cell = 'D2'; Value = '26.322.51'; Numeric = $false },
    @{ Cell = 'E2'; Value = '  -3.25%  '; Numeric = $false },
    @{ Cell = 'D3'; Value = '1.830.51'; Numeric = $false },
    @{ Cell = 'E3'; Value = '  -2.79%  '; Numeric = $false },
    @{ Cell = 'E4'; Value = '  +0.14%  '; Numeric = $false },
    @{ Cell = 'D5'; Value = '258.59'; Numeric = $true },
    @{ Cell = 'E5'; Value = '  -8.22%  '; Numeric = $false },
    @{ Cell = 'E6'; Value = '  +0.09%  '; Numeric = $false },
    @{ Cell = 'D7'; Value = '0.5199'; Numeric = $true },
    @{ Cell = 'E7'; Value = '  -1.99%  '; Numeric = $false },
    @{ Cell = 'D8'; Value = '0.3221'; Numeric = $true },
    @{ Cell = 'E8'; Value = '  -8.99%  '; Numeric = $false },
    @{ Cell = 'D9'; Value = '0.06719'; Numeric = $true },
    @{ Cell = 'E9'; Value = '  -4.74%  '; Numeric = $false },
    @{ Cell = 'D10'; Value = '18.65'; Numeric = $true },
    @{ Cell = 'E10'; Value = '  -8.76%  '; Numeric = $false },
    @{ Cell = 'D11'; Value = '0.7644'; Numeric = $true },
    @{ Cell = 'E11'; Value = '  -7.12%  '; Numeric = $false },
    @{ Cell = 'B12'; Value = 'TRON'; Numeric = $false },
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; Numeric = $false },
    @{ Cell = 'D12'; Value = '0.07672'; Numeric = $true },
    @{ Cell = 'E12'; Value = '  -1.94%  '; Numeric = $false },
    @{ Cell = 'B13'; Value = 'WrappedEther'; Numeric = $false },
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; Numeric = $false },
    @{ Cell = 'D13'; Value = '1.858.39'; Numeric = $false },
    @{ Cell = 'E13'; Value = '  -1.31%  '; Numeric = $false },
    @{ Cell = 'D14'; Value = '88.59'; Numeric = $true },
    @{ Cell = 'E14'; Value = '  -2.76%  '; Numeric = $false },
    @{ Cell = 'D15'; Value = '5.016'; Numeric = $true },
    @{ Cell = 'E15'; Value = '  -3.82%  '; Numeric = $false },
    @{ Cell = 'D16'; Value = '1.000'; Numeric = $true },
    @{ Cell = 'E16'; Value = '  +0.12%  '; Numeric = $false },
    @{ Cell = 'E17'; Value = '  -4.06%  '; Numeric = $false },
    @{ Cell = 'E18'; Value = '  +0.16%  '; Numeric = $false },
    @{ Cell = 'D19'; Value = '0.000007887'; Numeric = $true },
    @{ Cell = 'D20'; Value = '26.361.82'; Numeric = $false },
    @{ Cell = 'E20'; Value = '  -3.19%  '; Numeric = $false },
    @{ Cell = 'D21'; Value = '2.086.09'; Numeric = $false },
    @{ Cell = 'E21'; Value = '  -2.93%  '; Numeric = $false },
    @{ Cell = 'D22'; Value = '4.522'; Numeric = $true },
    @{ Cell = 'E22'; Value = '  -5.45%  '; Numeric = $false },
    @{ Cell = 'D23'; Value = '9.403'; Numeric = $true },
    @{ Cell = 'E23'; Value = '  -7.57%  '; Numeric = $false },
    @{ Cell = 'D24'; Value = '5.891'; Numeric = $true },
    @{ Cell = 'E24'; Value = '  -5.79%  '; Numeric = $false },
    @{ Cell = 'D25'; Value = '2.279'; Numeric = $true },
    @{ Cell = 'E25'; Value = '  -5.48%  '; Numeric = $false },
    @{ Cell = 'D26'; Value = '145.16'; Numeric = $true },
    @{ Cell = 'E26'; Value = '  -1.53%  '; Numeric = $false },
    @{ Cell = 'D27'; Value = '1.638'; Numeric = $true },
    @{ Cell = 'E27'; Value = '  -2.25%  '; Numeric = $false },
    @{ Cell = 'D28'; Value = '16.88'; Numeric = $true },
    @{ Cell = 'E28'; Value = '  -4.35%  '; Numeric = $false },
    @{ Cell = 'D29'; Value = '110.83'; Numeric = $true },
    @{ Cell = 'E29'; Value = '  -3.07%  '; Numeric = $false },
    @{ Cell = 'D30'; Value = '4.183'; Numeric = $true },
    @{ Cell = 'E30'; Value = '  -5.63%  '; Numeric = $false },
    @{ Cell = 'D31'; Value = '4.119'; Numeric = $true },
    @{ Cell = 'E31'; Value = '  -6.47%  '; Numeric = $false },
    @{ Cell = 'D32'; Value = '0.08720'; Numeric = $true },
    @{ Cell = 'E32'; Value = '  -2.77%  '; Numeric = $false },
    @{ Cell = 'D33'; Value = '0.04829'; Numeric = $true },
    @{ Cell = 'E33'; Value = '  -2.51%  '; Numeric = $false },
    @{ Cell = 'E34'; Value = '  -5.32%  '; Numeric = $false },
    @{ Cell = 'D35'; Value = '2.842'; Numeric = $true },
    @{ Cell = 'E35'; Value = '  -2.18%  '; Numeric = $false },
    @{ Cell = 'D36'; Value = '0.6802'; Numeric = $true },
    @{ Cell = 'E36'; Value = '  -9.38%  '; Numeric = $false },
    @{ Cell = 'D37'; Value = '3.086'; Numeric = $true },
    @{ Cell = 'E37'; Value = '  -6.93%  '; Numeric = $false },
    @{ Cell = 'E38'; Value = '  -6.13%  '; Numeric = $false },
    @{ Cell = 'D39'; Value = '2.207'; Numeric = $true },
    @{ Cell = 'E39'; Value = '  -9.16%  '; Numeric = $false },
    @{ Cell = 'D40'; Value = '0.4885'; Numeric = $true },
    @{ Cell = 'E40'; Value = '  -8.53%  '; Numeric = $false },
    @{ Cell = 'D41'; Value = '112.21'; Numeric = $true },
    @{ Cell = 'E41'; Value = '  -4.20%  '; Numeric = $false },
    @{ Cell = 'B42'; Value = 'TrustWalletToken'; Numeric = $false },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; Numeric = $false },
    @{ Cell = 'D42'; Value = '0.8878'; Numeric = $true },
    @{ Cell = 'E42'; Value = '  -8.80%  '; Numeric = $false },
    @{ Cell = 'B43'; Value = 'FraxShare'; Numeric = $false },
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; Numeric = $false },
    @{ Cell = 'D43'; Value = '6.108'; Numeric = $true },
    @{ Cell = 'E43'; Value = '  -3.68%  '; Numeric = $false },
    @{ Cell = 'D44'; Value = '0.9998'; Numeric = $true },
    @{ Cell = 'E44'; Value = '  +0.09%  '; Numeric = $false },
    @{ Cell = 'D45'; Value = '7.686'; Numeric = $true },
    @{ Cell = 'E45'; Value = '  -6.98%  '; Numeric = $false },
    @{ Cell = 'D46'; Value = '0.4187'; Numeric = $true },
    @{ Cell = 'E46'; Value = '  -9.54%  '; Numeric = $false },
    @{ Cell = 'D47'; Value = '0.1252'; Numeric = $true },
    @{ Cell = 'E47'; Value = '  -8.80%  '; Numeric = $false },
    @{ Cell = 'B48'; Value = 'EnergySwap'; Numeric = $false },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; Numeric = $false },
    @{ Cell = 'D48'; Value = '9.058'; Numeric = $true },
    @{ Cell = 'E48'; Value = '  -4.37%  '; Numeric = $false },
    @{ Cell = 'B49'; Value = 'Cronos'; Numeric = $false },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; Numeric = $false },
    @{ Cell = 'D49'; Value = '0.05874'; Numeric = $true },
    @{ Cell = 'E49'; Value = '  -1.45%  '; Numeric = $false },
    @{ Cell = 'D50'; Value = '35.27'; Numeric = $true },
    @{ Cell = 'E50'; Value = '  -4.17%  '; Numeric = $false },
    @{ Cell = 'D51'; Value = '59.18'; Numeric = $true },
    @{ Cell = 'E51'; Value = '  -4.60%  '; Numeric = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Numeric) {
        # Force text storage so strings like "1.000" or "18.65" keep their
        # exact literal form instead of becoming numeric values.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
